$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "23.464.28"
Set-TextValue $ws "E2" "  -1.42%  "
Set-TextValue $ws "D3" "1.647.13"
Set-TextValue $ws "E3" "  -0.41%  "
Set-TextValue $ws "E4" "  +0.05%  "
Set-TextValue $ws "D5" "1.001"
Set-TextValue $ws "E5" "  +0.06%  "
Set-TextValue $ws "D6" "298.82"
Set-TextValue $ws "E6" "  -1.86%  "
Set-TextValue $ws "D7" "0.3781"
Set-TextValue $ws "E7" "  -1.00%  "
Set-TextValue $ws "D8" "0.3538"
Set-TextValue $ws "E8" "  -2.04%  "
Set-TextValue $ws "D9" "50.03"
Set-TextValue $ws "E9" "  -2.24%  "
Set-TextValue $ws "D10" "0.08085"
Set-TextValue $ws "E10" "  -1.73%  "
Set-TextValue $ws "D11" "1.214"
Set-TextValue $ws "E11" "  -2.88%  "
Set-TextValue $ws "D12" "1.001"
Set-TextValue $ws "E12" "  +0.12%  "
Set-TextValue $ws "D13" "22.08"
Set-TextValue $ws "E13" "  -2.79%  "
Set-TextValue $ws "D14" "6.393"
Set-TextValue $ws "E14" "  -2.24%  "
Set-TextValue $ws "D15" "7.341"
Set-TextValue $ws "E15" "  -1.11%  "
Set-TextValue $ws "E16" "  -3.08%  "
Set-TextValue $ws "D17" "1.646.52"
Set-TextValue $ws "E17" "  +0.15%  "
Set-TextValue $ws "D18" "97.02"
Set-TextValue $ws "E18" "  -0.82%  "
Set-TextValue $ws "D19" "0.06955"
Set-TextValue $ws "E19" "  -0.31%  "
Set-TextValue $ws "D20" "6.758"
Set-TextValue $ws "E20" "  -0.29%  "
Set-TextValue $ws "E21" "  -2.25%  "
Set-TextValue $ws "D22" "1.000"
Set-TextValue $ws "E22" "  +0.00%  "
Set-TextValue $ws "D23" "12.44"
Set-TextValue $ws "E23" "  -2.17%  "
Set-TextValue $ws "D24" "23.466.46"
Set-TextValue $ws "D25" "2.496"
Set-TextValue $ws "E25" "  -1.80%  "
Set-TextValue $ws "D26" "2.890"
Set-TextValue $ws "E26" "  -6.74%  "
Set-TextValue $ws "D27" "20.87"
Set-TextValue $ws "E27" "  -2.12%  "
Set-TextValue $ws "D28" "152.03"
Set-TextValue $ws "E28" "  +0.62%  "
Set-TextValue $ws "D29" "5.197"
Set-TextValue $ws "E29" "  -0.94%  "
Set-TextValue $ws "D30" "132.79"
Set-TextValue $ws "E30" "  -1.44%  "
Set-TextValue $ws "D31" "1.828.18"
Set-TextValue $ws "E31" "  -0.13%  "
Set-TextValue $ws "D32" "6.943"
Set-TextValue $ws "E32" "  +0.56%  "
Set-TextValue $ws "D33" "2.148"
Set-TextValue $ws "E33" "  +1.23%  "
Set-TextValue $ws "D34" "11.47"
Set-TextValue $ws "E34" "  -3.88%  "
Set-TextValue $ws "D35" "0.9901"
Set-TextValue $ws "E35" "  -8.84%  "
Set-TextValue $ws "D36" "0.02715"
Set-TextValue $ws "E36" "  -4.85%  "
Set-TextValue $ws "D37" "0.08757"
Set-TextValue $ws "E37" "  -0.91%  "
Set-TextValue $ws "D38" "0.2438"
Set-TextValue $ws "E38" "  -3.39%  "
Set-TextValue $ws "D39" "5.934"
Set-TextValue $ws "E39" "  -3.66%  "
Set-TextValue $ws "D40" "12.96"
Set-TextValue $ws "E40" "  +0.15%  "
Set-TextValue $ws "D41" "0.06790"
Set-TextValue $ws "E41" "  -5.05%  "
Set-TextValue $ws "D42" "0.6885"
Set-TextValue $ws "E42" "  -2.91%  "
Set-TextValue $ws "B43" "TrustWalletToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D43" "1.296"
Set-TextValue $ws "E43" "  -3.48%  "
Set-TextValue $ws "B44" "EnergySwap"
Set-TextValue $ws "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D44" "15.69"
Set-TextValue $ws "E44" "  -1.17%  "
Set-TextValue $ws "D46" "0.6363"
Set-TextValue $ws "E46" "  -3.07%  "
Set-TextValue $ws "D47" "2.251"
Set-TextValue $ws "E47" "  -3.75%  "
Set-TextValue $ws "D48" "3.908"
Set-TextValue $ws "E48" "  -1.45%  "
Set-TextValue $ws "D49" "0.07724"
Set-TextValue $ws "E49" "  -3.36%  "
Set-TextValue $ws "D50" "127.70"
Set-TextValue $ws "E50" "  -0.93%  "
Set-TextValue $ws "D51" "1.149"
Set-TextValue $ws "E51" "  -4.05%  "
